# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de) the localization run has now
# produced a handback: the Status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", a "Latest Target File" (F) and
# "Latest Handback File" (G) column get populated (mirroring the source
# markdown file and the handed-off xlf file, each as a hyperlink just like
# the existing ones), and "Latest Handback DateTime" (H) gets a real
# timestamp instead of the 0001-01-01 00:00:00 placeholder.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status

$wsZh.Range("H2").Value = "2016-03-20 00:48:17"
$wsZh.Range("H3").Value = "2016-03-20 00:48:17"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/2bc49e7c-1f69-4273-ba0d-714a75eb899d.md",
    "",
    "",
    "2bc49e7c-1f69-4273-ba0d-714a75eb899d.md")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2b44fb760139019c3e7b68353811d834b3af4dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.zh-cn.xlf",
    "",
    "",
    "2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.zh-cn.xlf")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/e7faba28-3992-4225-bdb9-c5f7d617bd3e.md",
    "",
    "",
    "e7faba28-3992-4225-bdb9-c5f7d617bd3e.md")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2b44fb760139019c3e7b68353811d834b3af4dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.zh-cn.xlf",
    "",
    "",
    "e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.zh-cn.xlf")

$wsZh.Range("F2:G3").Style = "HyperLink"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status

$wsDe.Range("H2").Value = "2016-03-20 00:48:23"
$wsDe.Range("H3").Value = "2016-03-20 00:48:23"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/2bc49e7c-1f69-4273-ba0d-714a75eb899d.md",
    "",
    "",
    "2bc49e7c-1f69-4273-ba0d-714a75eb899d.md")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e10cc10f777b075f208c7e22e5aa2f5fdeaaf812/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.de-de.xlf",
    "",
    "",
    "2bc49e7c-1f69-4273-ba0d-714a75eb899d.8a4e2a08faba8868e51f71545f50f9b9092f20f5.de-de.xlf")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0b7a7e77138c1f6f859b0c9a076eff12bef6ecf3/e2e/e7faba28-3992-4225-bdb9-c5f7d617bd3e.md",
    "",
    "",
    "e7faba28-3992-4225-bdb9-c5f7d617bd3e.md")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e10cc10f777b075f208c7e22e5aa2f5fdeaaf812/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.de-de.xlf",
    "",
    "",
    "e7faba28-3992-4225-bdb9-c5f7d617bd3e.22b7535328d6b6d3b62c784bcce841610a6765dd.de-de.xlf")

$wsDe.Range("F2:G3").Style = "HyperLink"

Write-Output "Handback report generated for zh-cn and de-de sheets."
